# Insert two new data rows (252 and 253) into the "Vega Monumental Concepción -
# Zanahoria" weekly sheet, pushing all the existing rows from 252 downward by
# two positions (old row 252 -> new row 254, ... old row 358 -> new row 360).
# The sheet grows from A1:R358 to A1:R360.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from row 252 onward down by two rows.
$ws.Range("A252:A253").EntireRow.Insert()

# New row 252: "Primera" quality entry for the newly-reported week.
$ws.Cells.Item(252, 1).Value  = 11
$ws.Cells.Item(252, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(252, 3).Value  = "Bíobío"
$ws.Cells.Item(252, 4).Value  = 45006
$ws.Cells.Item(252, 5).Value  = 8
$ws.Cells.Item(252, 6).Value  = 100114013
$ws.Cells.Item(252, 7).Value  = "Zanahoria"
$ws.Cells.Item(252, 8).Value  = "Sin especificar"
$ws.Cells.Item(252, 9).Value  = "Primera"
$ws.Cells.Item(252, 10).Value = 600
$ws.Cells.Item(252, 11).Value = 6000
$ws.Cells.Item(252, 12).Value = 6500
$ws.Cells.Item(252, 13).Value = 6250
$ws.Cells.Item(252, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(252, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(252, 16).Value = 312
$ws.Cells.Item(252, 17).Value = 20
$ws.Cells.Item(252, 18).Value = "Hortaliza"

# New row 253: "Segunda" quality entry for the same newly-reported week.
$ws.Cells.Item(253, 1).Value  = 11
$ws.Cells.Item(253, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(253, 3).Value  = "Bíobío"
$ws.Cells.Item(253, 4).Value  = 45006
$ws.Cells.Item(253, 5).Value  = 8
$ws.Cells.Item(253, 6).Value  = 100114013
$ws.Cells.Item(253, 7).Value  = "Zanahoria"
$ws.Cells.Item(253, 8).Value  = "Sin especificar"
$ws.Cells.Item(253, 9).Value  = "Segunda"
$ws.Cells.Item(253, 10).Value = 300
$ws.Cells.Item(253, 11).Value = 5000
$ws.Cells.Item(253, 12).Value = 5000
$ws.Cells.Item(253, 13).Value = 5000
$ws.Cells.Item(253, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(253, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(253, 16).Value = 250
$ws.Cells.Item(253, 17).Value = 20
$ws.Cells.Item(253, 18).Value = "Hortaliza"
